$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation prompt
$excel.DisplayAlerts = $false

# Update the regression weights/bias on the remaining sheet (Layer0)
$ws = $wb.Worksheets.Item("Layer0")
$ws.Range("B2").Value = -0.0772216348024467
$ws.Range("C2").Value = -0.4719153656784933
$ws.Range("B3").Value = -0.3244479731105778
$ws.Range("C3").Value = -0.339491093008931
$ws.Range("B4").Value = -0.73218848422015
$ws.Range("C4").Value = -0.1255680597707436

# Remove the Layer1 sheet entirely
$ws2 = $wb.Worksheets.Item("Layer1")
[void]$ws2.Delete()

$excel.DisplayAlerts = $true
